$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$timestamp = "2026-02-20T23:33:43.608782"

# Row 2
$ws.Range("I2").Value = 1
$ws.Range("J2").Value = 0.00006944444444444444
$ws.Range("K2").Value = 469
$ws.Range("L2").Value = 0.000938
$ws.Range("M2").Value = $timestamp

# Row 3
$ws.Range("K3").Value = 495
$ws.Range("L3").Value = 0.00495
$ws.Range("M3").Value = $timestamp

# Row 4
$ws.Range("I4").Value = 0
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 0
$ws.Range("L4").Value = 0
$ws.Range("M4").Value = $timestamp

# Row 5
$ws.Range("I5").Value = 0
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 0
$ws.Range("L5").Value = 0
$ws.Range("M5").Value = $timestamp

# Row 6
$ws.Range("M6").Value = $timestamp

# Row 7
$ws.Range("M7").Value = $timestamp

# Row 8
$ws.Range("M8").Value = $timestamp

# Row 9
$ws.Range("I9").Value = 0
$ws.Range("J9").Value = 0
$ws.Range("K9").Value = 0
$ws.Range("L9").Value = 0
$ws.Range("M9").Value = $timestamp

# Row 10
$ws.Range("I10").Value = 0
$ws.Range("J10").Value = 0
$ws.Range("K10").Value = 0
$ws.Range("L10").Value = 0
$ws.Range("M10").Value = $timestamp

# Row 11
$ws.Range("I11").Value = 0
$ws.Range("J11").Value = 0
$ws.Range("K11").Value = 0
$ws.Range("L11").Value = 0
$ws.Range("M11").Value = $timestamp

# Row 12
$ws.Range("I12").Value = 0
$ws.Range("J12").Value = 0
$ws.Range("K12").Value = 0
$ws.Range("L12").Value = 0
$ws.Range("M12").Value = $timestamp

# Row 13
$ws.Range("M13").Value = $timestamp

# Row 14
$ws.Range("I14").Value = 0
$ws.Range("J14").Value = 0
$ws.Range("K14").Value = 0
$ws.Range("L14").Value = 0
$ws.Range("M14").Value = $timestamp

$wb.Save()
